$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 3500
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 3500
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 3500
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -4086

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5000
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 101782.3
$ws.Range("I101").Value = 1834
$ws.Range("J101").Value = 334995
$ws.Range("K101").Value = 5502
$ws.Range("L101").Value = 1004985
$ws.Range("M101").Value = -3880
$ws.Range("N101").Value = -1008229

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 45781.11
$ws.Range("J123").Value = 45781.11
$ws.Range("L123").Value = 45781.11
$ws.Range("N123").Value = -55581.11

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 33956.668
$ws.Range("J128").Value = 33956.668
$ws.Range("L128").Value = 33956.668
$ws.Range("N128").Value = -43916.668

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 992.5484
$ws.Range("I135").Value = 1034.7307
$ws.Range("J135").Value = 773.2
$ws.Range("K135").Value = 9312.576300000001
$ws.Range("L135").Value = 6958.8
$ws.Range("M135").Value = -6777.576300000001
$ws.Range("N135").Value = -12028.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1575.909
$ws.Range("I141").Value = 1575.909
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4727.727000000001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 452.2729999999992
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5083.86
$ws.Range("I32").Value = 4881.4897
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 4881.4897
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -4594.4897
$ws.Range("N32").Value = -15574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1043.9
$ws.Range("I45").Value = 929.0769
$ws.Range("J45").Value = 1257.1428
$ws.Range("K45").Value = 929.0769
$ws.Range("L45").Value = 1257.1428
$ws.Range("M45").Value = -552.0769
$ws.Range("N45").Value = -2011.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 8600
$ws.Range("I57").Value = 8600
$ws.Range("K57").Value = 8600
$ws.Range("M57").Value = -8116

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1360.4348
$ws.Range("I74").Value = 1399.4736
$ws.Range("J74").Value = 1175
$ws.Range("K74").Value = 1399.4736
$ws.Range("L74").Value = 1175
$ws.Range("M74").Value = -525.4736
$ws.Range("N74").Value = -2923

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1360.4348
$ws.Range("I77").Value = 1399.4736
$ws.Range("J77").Value = 1175
$ws.Range("K77").Value = 6997.368
$ws.Range("L77").Value = 5875
$ws.Range("M77").Value = -2629.368
$ws.Range("N77").Value = -14611

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 23359.8
$ws.Range("I113").Value = 23359.8
$ws.Range("K113").Value = 23359.8
$ws.Range("M113").Value = -21189.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2419.4827
$ws.Range("I134").Value = 1932.037
$ws.Range("J134").Value = 9000
$ws.Range("K134").Value = 5796.111
$ws.Range("L134").Value = 27000
$ws.Range("M134").Value = -3261.111
$ws.Range("N134").Value = -32070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2527.0984
$ws.Range("I58").Value = 894.44446
$ws.Range("J58").Value = 4878.12
$ws.Range("K58").Value = 894.44446
$ws.Range("L58").Value = 4878.12
$ws.Range("M58").Value = -691.44446
$ws.Range("N58").Value = -5284.12

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3083.1333
$ws.Range("I99").Value = 2116.6667
$ws.Range("J99").Value = 4532.8335
$ws.Range("K99").Value = 2116.6667
$ws.Range("L99").Value = 4532.8335
$ws.Range("M99").Value = -618.6667000000002
$ws.Range("N99").Value = -7528.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3083.1333
$ws.Range("I126").Value = 2116.6667
$ws.Range("J126").Value = 4532.8335
$ws.Range("K126").Value = 6350.000100000001
$ws.Range("L126").Value = 13598.5005
$ws.Range("M126").Value = -3880.000100000001
$ws.Range("N126").Value = -18538.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 20834778
$ws.Range("I134").Value = 1250.6316
$ws.Range("J134").Value = 100002184
$ws.Range("K134").Value = 3751.8948
$ws.Range("L134").Value = 300006552
$ws.Range("M134").Value = -1216.8948
$ws.Range("N134").Value = -300011622

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2527.0984
$ws.Range("I136").Value = 894.44446
$ws.Range("J136").Value = 4878.12
$ws.Range("K136").Value = 2683.33338
$ws.Range("L136").Value = 14634.36
$ws.Range("M136").Value = -133.33338
$ws.Range("N136").Value = -19734.36

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4059.4285
$ws.Range("I3").Value = 1534.4615
$ws.Range("J3").Value = 8162.5
$ws.Range("K3").Value = 4603.3845
$ws.Range("L3").Value = 24487.5
$ws.Range("M3").Value = -4491.3845
$ws.Range("N3").Value = -24711.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 143.28572
$ws.Range("I40").Value = 145
$ws.Range("J40").Value = 141
$ws.Range("K40").Value = 580
$ws.Range("L40").Value = 564
$ws.Range("M40").Value = -511
$ws.Range("N40").Value = -702

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 348
$ws.Range("I97").Value = 340.7143
$ws.Range("K97").Value = 1022.1429
$ws.Range("M97").Value = -526.1428999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 38969.58
$ws.Range("J123").Value = 38969.58
$ws.Range("L123").Value = 38969.58
$ws.Range("N123").Value = -43869.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2487.2368
$ws.Range("I132").Value = 2358.4583
$ws.Range("J132").Value = 2708
$ws.Range("K132").Value = 7075.374899999999
$ws.Range("L132").Value = 8124
$ws.Range("M132").Value = -4545.374899999999
$ws.Range("N132").Value = -13184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 11999.556
$ws.Range("J136").Value = 11999.556
$ws.Range("L136").Value = 35998.66800000001
$ws.Range("N136").Value = -41098.66800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 47266.668
$ws.Range("J138").Value = 47266.668
$ws.Range("L138").Value = 47266.668
$ws.Range("N138").Value = -57546.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1065.9286
$ws.Range("I22").Value = 307.2857
$ws.Range("K22").Value = 307.2857
$ws.Range("M22").Value = -12.28570000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1065.9286
$ws.Range("I27").Value = 307.2857
$ws.Range("K27").Value = 307.2857
$ws.Range("M27").Value = -200.2857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 28999.5
$ws.Range("J115").Value = 28999.5
$ws.Range("L115").Value = 28999.5
$ws.Range("N115").Value = -31349.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4057.9783
$ws.Range("I136").Value = 2248.8215
$ws.Range("J136").Value = 6872.222
$ws.Range("K136").Value = 6746.4645
$ws.Range("L136").Value = 20616.666
$ws.Range("M136").Value = -4196.4645
$ws.Range("N136").Value = -25716.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2003.1538
$ws.Range("I81").Value = 1949
$ws.Range("J81").Value = 2125
$ws.Range("K81").Value = 3898
$ws.Range("L81").Value = 4250
$ws.Range("M81").Value = -2837
$ws.Range("N81").Value = -6372

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2003.1538
$ws.Range("I84").Value = 1949
$ws.Range("J84").Value = 2125
$ws.Range("K84").Value = 19490
$ws.Range("L84").Value = 21250
$ws.Range("M84").Value = -14186
$ws.Range("N84").Value = -31858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 1031200
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 1031200
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 1031200
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -1041020

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 762.4545000000001
$ws.Range("I132").Value = 659.4878
$ws.Range("J132").Value = 2169.6667
$ws.Range("K132").Value = 1978.4634
$ws.Range("L132").Value = 6509.000100000001
$ws.Range("M132").Value = 551.5365999999999
$ws.Range("N132").Value = -11569.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1791.4736
$ws.Range("I136").Value = 2231.2856
$ws.Range("J136").Value = 560
$ws.Range("K136").Value = 6693.8568
$ws.Range("L136").Value = 1680
$ws.Range("M136").Value = -4143.8568
$ws.Range("N136").Value = -6780
